$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2-17: Price (D) and Volume (E) columns
$ws.Range('D2').Value = '26.274.20'
$ws.Range('E2').Value = '  +0.58%  '
$ws.Range('D3').Value = '1.663.09'
$ws.Range('E3').Value = '  +0.44%  '
$ws.Range('D4').Value = "'1.011"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.83%  '
$ws.Range('D5').Value = "'218.73"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.31%  '
$ws.Range('D6').Value = "'0.5301"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.25%  '
$ws.Range('D7').Value = "'1.011"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.79%  '
$ws.Range('D8').Value = "'0.2634"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.07%  '
$ws.Range('D9').Value = "'0.06357"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.33%  '
$ws.Range('D10').Value = "'20.52"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.58%  '
$ws.Range('D11').Value = "'0.07855"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.93%  '
$ws.Range('D12').Value = "'4.561"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.45%  '
$ws.Range('D13').Value = '1.663.81'
$ws.Range('E13').Value = '  +0.42%  '
$ws.Range('D14').Value = '1.891.92'
$ws.Range('E14').Value = '  +0.46%  '
$ws.Range('D15').Value = "'0.5521"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.90%  '
$ws.Range('D16').Value = '0.0₅8170'
$ws.Range('E16').Value = '  +0.06%  '
$ws.Range('D17').Value = "'65.62"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.22%  '

# Rows 18-51: Coin (B), Link (C), Price (D), Volume (E) all shift due to new WrappedBTC row insertion
$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').Value = '26.287.61'
$ws.Range('E18').Value = '  +0.58%  '
$ws.Range('B19').Value = 'Dai'
$ws.Range('C19').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D19').Value = "'1.010"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.79%  '
$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D20').Value = "'4.669"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.33%  '
$ws.Range('B21').Value = 'BitcoinCash'
$ws.Range('C21').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D21').Value = "'192.74"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.06%  '
$ws.Range('B22').Value = 'Avalanche'
$ws.Range('C22').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D22').Value = "'10.21"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.17%  '
$ws.Range('B23').Value = 'Chainlink'
$ws.Range('C23').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D23').Value = "'6.028"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.10%  '
$ws.Range('B24').Value = 'BinanceUSD'
$ws.Range('C24').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D24').Value = "'1.012"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.82%  '
$ws.Range('B25').Value = 'Monero'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D25').Value = "'144.28"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.72%  '
$ws.Range('B26').Value = 'Stellar'
$ws.Range('C26').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D26').Value = "'0.1226"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.04%  '
$ws.Range('B27').Value = 'Cosmos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D27').Value = "'7.192"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.04%  '
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').Value = "'16.06"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.86%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').Value = "'1.478"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.71%  '
$ws.Range('B30').Value = 'Hedera'
$ws.Range('C30').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D30').Value = "'0.05918"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.31%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').Value = "'1.283"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.38%  '
$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D32').Value = "'3.587"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.14%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').Value = "'3.275"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.90%  '
$ws.Range('B34').Value = 'LidoDAOToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D34').Value = "'1.614"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.56%  '
$ws.Range('B35').Value = 'MXToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D35').Value = "'2.827"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.18%  '
$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D36').Value = "'0.9583"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.97%  '
$ws.Range('B37').Value = 'HuobiToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D37').Value = "'2.426"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.59%  '
$ws.Range('B38').Value = 'ImmutableX'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D38').Value = "'0.5786"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.08%  '
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').Value = "'0.01601"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.56%  '
$ws.Range('B40').Value = 'TrustWalletToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D40').Value = "'0.8651"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.07%  '
$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D41').Value = "'5.854"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.66%  '
$ws.Range('B42').Value = 'PaxDollar'
$ws.Range('C42').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D42').Value = "'1.010"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.75%  '
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').Value = '1.045.30'
$ws.Range('E43').Value = '  +2.29%  '
$ws.Range('B44').Value = 'Quant'
$ws.Range('C44').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D44').Value = "'103.97"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.47%  '
$ws.Range('B45').Value = 'RocketPoolETH'
$ws.Range('C45').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D45').Value = '1.804.91'
$ws.Range('E45').Value = '  +0.34%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').Value = "'57.35"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.33%  '
$ws.Range('B47').Value = 'BabyDogeCoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D47').Value = '0.0₈106'
$ws.Range('E47').Value = '  -4.96%  '
$ws.Range('B48').Value = 'Frax'
$ws.Range('C48').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D48').Value = "'1.011"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.35%  '
$ws.Range('B49').Value = 'Mantle'
$ws.Range('C49').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D49').Value = "'0.4382"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.16%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = "'8.023"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.65%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').Value = "'0.05164"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.27%  '

Write-Host "Update complete"